$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 203.7816646666667
$ws.Cells.Item(2, 8).Value = 611.344994
$ws.Cells.Item(2, 9).Value = 0.6667327591988204
$ws.Cells.Item(2, 10).Value = 0.6667327591988205
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 1.508394
$ws.Cells.Item(2, 14).Value = 4.525182
$ws.Cells.Item(2, 15).Value = 0.0276475339394655
$ws.Cells.Item(2, 16).Value = 0.0276475339394655
$ws.Cells.Item(2, 17).Value = 307.383040293212
$ws.Cells.Item(2, 18).Value = 2766.447362638908
$ws.Cells.Item(2, 19).Value = 0.01843351658850287
$ws.Cells.Item(2, 20).Value = 0.01843351658850287

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 203.7816646666667
$ws.Cells.Item(3, 8).Value = 611.344994
$ws.Cells.Item(3, 9).Value = 0.6667327591988204
$ws.Cells.Item(3, 10).Value = 0.6667327591988205
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 8.961352
$ws.Cells.Item(3, 14).Value = 26.884056
$ws.Cells.Item(3, 15).Value = 0.1642536920482958
$ws.Cells.Item(3, 16).Value = 0.1642536920482958
$ws.Cells.Item(3, 17).Value = 1826.159228223963
$ws.Cells.Item(3, 18).Value = 16435.43305401567
$ws.Cells.Item(3, 19).Value = 0.1095133173079536
$ws.Cells.Item(3, 20).Value = 0.1095133173079536

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 203.7816646666667
$ws.Cells.Item(4, 8).Value = 611.344994
$ws.Cells.Item(4, 9).Value = 0.6667327591988204
$ws.Cells.Item(4, 10).Value = 0.6667327591988205
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 44.08824833333333
$ws.Cells.Item(4, 14).Value = 132.264745
$ws.Cells.Item(4, 15).Value = 0.8080987740122387
$ws.Cells.Item(4, 16).Value = 0.8080987740122386
$ws.Cells.Item(4, 17).Value = 8984.376637604058
$ws.Cells.Item(4, 18).Value = 80859.38973843654
$ws.Cells.Item(4, 19).Value = 0.5387859253023639
$ws.Cells.Item(4, 20).Value = 0.5387859253023639

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 63.14058933333333
$ws.Cells.Item(5, 8).Value = 189.421768
$ws.Cells.Item(5, 9).Value = 0.2065833519051582
$ws.Cells.Item(5, 10).Value = 0.2065833519051582
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 1.508394
$ws.Cells.Item(5, 14).Value = 4.525182
$ws.Cells.Item(5, 15).Value = 0.0276475339394655
$ws.Cells.Item(5, 16).Value = 0.0276475339394655
$ws.Cells.Item(5, 17).Value = 95.240886106864
$ws.Cells.Item(5, 18).Value = 857.1679749617759
$ws.Cells.Item(5, 19).Value = 0.005711520233126407
$ws.Cells.Item(5, 20).Value = 0.005711520233126407

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 63.14058933333333
$ws.Cells.Item(6, 8).Value = 189.421768
$ws.Cells.Item(6, 9).Value = 0.2065833519051582
$ws.Cells.Item(6, 10).Value = 0.2065833519051582
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 8.961352
$ws.Cells.Item(6, 14).Value = 26.884056
$ws.Cells.Item(6, 15).Value = 0.1642536920482958
$ws.Cells.Item(6, 16).Value = 0.1642536920482958
$ws.Cells.Item(6, 17).Value = 565.8250465034453
$ws.Cells.Item(6, 18).Value = 5092.425418531007
$ws.Cells.Item(6, 19).Value = 0.03393207826613457
$ws.Cells.Item(6, 20).Value = 0.03393207826613457

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 63.14058933333333
$ws.Cells.Item(7, 8).Value = 189.421768
$ws.Cells.Item(7, 9).Value = 0.2065833519051582
$ws.Cells.Item(7, 10).Value = 0.2065833519051582
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 44.08824833333333
$ws.Cells.Item(7, 14).Value = 132.264745
$ws.Cells.Item(7, 15).Value = 0.8080987740122387
$ws.Cells.Item(7, 16).Value = 0.8080987740122386
$ws.Cells.Item(7, 17).Value = 2783.757982441018
$ws.Cells.Item(7, 18).Value = 25053.82184196916
$ws.Cells.Item(7, 19).Value = 0.1669397534058972
$ws.Cells.Item(7, 20).Value = 0.1669397534058972

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 38.719942
$ws.Cells.Item(8, 8).Value = 116.159826
$ws.Cells.Item(8, 9).Value = 0.1266838888960214
$ws.Cells.Item(8, 10).Value = 0.1266838888960214
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 1.508394
$ws.Cells.Item(8, 14).Value = 4.525182
$ws.Cells.Item(8, 15).Value = 0.0276475339394655
$ws.Cells.Item(8, 16).Value = 0.0276475339394655
$ws.Cells.Item(8, 17).Value = 58.404928193148
$ws.Cells.Item(8, 18).Value = 525.644353738332
$ws.Cells.Item(8, 19).Value = 0.003502497117836228
$ws.Cells.Item(8, 20).Value = 0.003502497117836229

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 38.719942
$ws.Cells.Item(9, 8).Value = 116.159826
$ws.Cells.Item(9, 9).Value = 0.1266838888960214
$ws.Cells.Item(9, 10).Value = 0.1266838888960214
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 8.961352
$ws.Cells.Item(9, 14).Value = 26.884056
$ws.Cells.Item(9, 15).Value = 0.1642536920482958
$ws.Cells.Item(9, 16).Value = 0.1642536920482958
$ws.Cells.Item(9, 17).Value = 346.983029681584
$ws.Cells.Item(9, 18).Value = 3122.847267134256
$ws.Cells.Item(9, 19).Value = 0.02080829647420761
$ws.Cells.Item(9, 20).Value = 0.02080829647420762

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 38.719942
$ws.Cells.Item(10, 8).Value = 116.159826
$ws.Cells.Item(10, 9).Value = 0.1266838888960214
$ws.Cells.Item(10, 10).Value = 0.1266838888960214
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 44.08824833333333
$ws.Cells.Item(10, 14).Value = 132.264745
$ws.Cells.Item(10, 15).Value = 0.8080987740122387
$ws.Cells.Item(10, 16).Value = 0.8080987740122386
$ws.Cells.Item(10, 17).Value = 1707.094418348263
$ws.Cells.Item(10, 18).Value = 15363.84976513437
$ws.Cells.Item(10, 19).Value = 0.1023730953039775
$ws.Cells.Item(10, 20).Value = 0.1023730953039776
